$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64: the phone number was previously stored as text; normalize it to a
# genuine number (matches how the payments app's next sync re-writes the row).
$ws.Range("A64").Value = 51616191

# Row 65: new payment for the same phone number, 51616191 (Cash).
# The phone column is written as text for freshly-appended rows, so enter it
# with a leading apostrophe (forces text) and then drop the resulting
# "quote prefix" formatting so the cell ends up plain text with no special
# number format applied.
$ws.Range("A65").Value = "'51616191"
$ws.Range("A65").ClearFormats()

$ws.Range("B65").Value = ""
$ws.Range("C65").Value = "Cash"
$ws.Range("D65").Value = "2025-08-20T08:04:39"
$ws.Range("E65").Value = 120
$ws.Range("F65").Value = ""
$ws.Range("G65").Value = 115
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 100
$ws.Range("J65").Value = 5
